# Actualización automática 2025-08-22 13:35:10
#
# Updates the monthly sales figures for advisor "ALMEIDA CUATIN JHONATHANN
# CARLOS": the INODOROS sale recorded for client MANCHENO PINO HERVIN
# SANTIAGO increased from 355.5 to 444.6, which ripples through the
# per-group totals (sheet 1), the monthly sales sheet (sheet 2) and the
# monthly compliance sheet (sheet 3).

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" -------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M15").Value = 1895.34
$wsGrupo.Range("H19").Value = 444.6
$wsGrupo.Range("M19").Value = 2021.43
$wsGrupo.Range("O24").Value = 1599.58

# --- Sheet "VENTA MENSUAL" -----------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F15").Value = 3799.99
$wsMensual.Range("F19").Value = 3172.07
$wsMensual.Range("F24").Value = 1599.58
$wsMensual.Range("F34").Value = 16473.46

# --- Sheet "CUMPLIMIENTO MENSUAL" ----------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumpl.Range("D7").Value = 444.6
$wsCumpl.Range("E7").Value = 1155.4
$wsCumpl.Range("F7").Value = 0.277875

$wsCumpl.Range("D16").Value = 7261.32
$wsCumpl.Range("E16").Value = 14611.78
$wsCumpl.Range("F16").Value = 0.3319748915334361

$wsCumpl.Range("D18").Value = 2001.17
$wsCumpl.Range("E18").Value = -401.1700000000001
$wsCumpl.Range("F18").Value = 1.25073125

$wsCumpl.Range("D19").Value = 16991.67
$wsCumpl.Range("E19").Value = 15117.61107555788
$wsCumpl.Range("F19").Value = 0.5291825114369921
